# TestDataSheet.xlsx maintenance edit:
#   - "Bangalore" renamed to "Bengaluru" everywhere it appears
#     (FlightBookingData!B3, HotelBookingData!A2, BookingDateToastMessageData!A2)
#   - Numeric "count" columns reformatted with an explicit integer
#     number format (Formatter) instead of the default General format
#   - Cursor/selection state left the way the author's last save left it

$wb = $excel.ActiveWorkbook

$flight  = $wb.Worksheets.Item("FlightBookingData")
$hotel   = $wb.Worksheets.Item("HotelBookingData")
$blog    = $wb.Worksheets.Item("TravelBlogData")
$booking = $wb.Worksheets.Item("BookingDateToastMessageData")
$guest   = $wb.Worksheets.Item("GuestLimitToastMessageData")

# --- Rename "Bangalore" -> "Bengaluru" -----------------------------------
$flight.Range("B3").Value = "Bengaluru"
$hotel.Range("A2").Value = "Bengaluru"
$booking.Range("A2").Value = "Bengaluru"

# --- Apply an integer Formatter to the numeric columns -------------------
$flight.Range("E2:G3").NumberFormat = "0"
$hotel.Range("F2:G3").NumberFormat = "0"
$guest.Range("F2").NumberFormat = "0"

# --- Restore the selection / active-cell state for each sheet ------------
# (touched in this order so the workbook ends up with GuestLimitToastMessageData
# as the active/selected tab, matching the saved file)
[void]$flight.Range("B4").Select()
[void]$hotel.Range("B7").Select()
[void]$booking.Range("A2").Select()
[void]$guest.Range("F2").Select()
